# "Bugifixed QoQ Visualizations and a typo in the evaluation objects"
#
# The evaluation series previously included four extra leading rows
# (year-end dates 1983-12-31 .. 1986-12-31, i.e. old rows 2-5) that
# don't belong in the series. Remove those four rows and shift the
# remaining data up, so the series now starts at 1987-12-31 (old row 6,
# which becomes the new row 2), and the sheet ends at row 39 instead of 43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:B5").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
